$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for the Price/Volume columns so numeric-looking
# strings (e.g. "616.23") are not coerced to floating point numbers,
# which would lose the exact textual representation used in the sheet.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '69.957.79'
$ws.Range("E2").Value = '  +0.20%  '
$ws.Range("D3").Value = '3.784.87'
$ws.Range("E3").Value = '  +3.74%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '616.23'
$ws.Range("E5").Value = '  +4.03%  '
$ws.Range("D6").Value = '178.33'
$ws.Range("E6").Value = '  -4.05%  '
$ws.Range("D7").Value = '3.784.67'
$ws.Range("E7").Value = '  +3.67%  '
$ws.Range("E8").Value = '  +0.11%  '
$ws.Range("D9").Value = '0.536'
$ws.Range("E9").Value = '  +0.14%  '
$ws.Range("D10").Value = '0.169'
$ws.Range("E10").Value = '  +5.07%  '
$ws.Range("D11").Value = '6.33'
$ws.Range("E11").Value = '  -3.11%  '
$ws.Range("D12").Value = '0.493'
$ws.Range("E12").Value = '  -0.60%  '
$ws.Range("D13").Value = '41.09'
$ws.Range("E13").Value = '  +4.64%  '
$ws.Range("D14").Value = '0.0000256'
$ws.Range("E14").Value = '  +1.17%  '
$ws.Range("D15").Value = '4.412.84'
$ws.Range("E15").Value = '  +3.97%  '
$ws.Range("D16").Value = '3.775.40'
$ws.Range("E16").Value = '  +3.92%  '
$ws.Range("D17").Value = '69.975.40'
$ws.Range("E17").Value = '  +0.18%  '
$ws.Range("D18").Value = '0.124'
$ws.Range("E18").Value = '  +0.00%  '
$ws.Range("D19").Value = '7.60'
$ws.Range("E19").Value = '  +0.62%  '
$ws.Range("D20").Value = '515.50'
$ws.Range("E20").Value = '  +1.24%  '
$ws.Range("D21").Value = '16.60'
$ws.Range("E21").Value = '  -3.19%  '
$ws.Range("D22").Value = '9.61'
$ws.Range("E22").Value = '  +2.82%  '
$ws.Range("D23").Value = '0.727'
$ws.Range("E23").Value = '  -2.44%  '
$ws.Range("E24").Value = '  +5.36%  '
$ws.Range("D25").Value = '88.21'
$ws.Range("E25").Value = '  +0.45%  '
$ws.Range("D26").Value = '13.36'
$ws.Range("E26").Value = '  -1.16%  '
$ws.Range("D27").Value = '11.10'
$ws.Range("E27").Value = '  +2.21%  '
$ws.Range("D28").Value = '0.0000135'
$ws.Range("E28").Value = '  +25.11%  '
$ws.Range("E29").Value = '  +0.09%  '
$ws.Range("D30").Value = '2.50'
$ws.Range("E30").Value = '  -1.73%  '
$ws.Range("D31").Value = '2.84'
$ws.Range("E31").Value = '  +3.45%  '
$ws.Range("D32").Value = '7.83'
$ws.Range("E32").Value = '  -4.16%  '
$ws.Range("D33").Value = '31.64'
$ws.Range("E33").Value = '  -3.09%  '
$ws.Range("E34").Value = '  -1.74%  '
$ws.Range("D35").Value = '0.998'
$ws.Range("E35").Value = '  -0.03%  '
$ws.Range("D36").Value = '6.23'
$ws.Range("E36").Value = '  +1.50%  '
$ws.Range("D37").Value = '1.05'
$ws.Range("E37").Value = '  +2.20%  '
$ws.Range("E38").Value = '  +1.43%  '
$ws.Range("D39").Value = '2.18'
$ws.Range("E39").Value = '  +3.05%  '
$ws.Range("D40").Value = '0.133'
$ws.Range("E40").Value = '  +3.91%  '
$ws.Range("E41").Value = '  +0.80%  '
$ws.Range("D42").Value = '44.43'
$ws.Range("E42").Value = '  -4.56%  '
$ws.Range("D43").Value = '8.80'
$ws.Range("E43").Value = '  -0.46%  '
$ws.Range("D44").Value = '424.03'
$ws.Range("E44").Value = '  +4.90%  '
$ws.Range("D45").Value = '3.064.80'
$ws.Range("E45").Value = '  -2.92%  '
$ws.Range("D46").Value = '2.75'
$ws.Range("E46").Value = '  -1.23%  '
$ws.Range("D47").Value = '0.0364'
$ws.Range("E47").Value = '  -0.69%  '
$ws.Range("D48").Value = '27.72'
$ws.Range("E48").Value = '  -0.29%  '
$ws.Range("D49").Value = '2.50'
$ws.Range("E49").Value = '  +2.13%  '
$ws.Range("E50").Value = '  -0.08%  '
$ws.Range("D51").Value = '135.59'
$ws.Range("E51").Value = '  -0.94%  '

# Restore number format / style so no stray formatting is left behind.
$dataRange.NumberFormat = "General"
$dataRange.Style = "Normal"

